# Auto-generated edit script: update crypto price/volume table per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.371.34"
$ws.Range("E2").Value = "  +5.23%  "
$ws.Range("D3").Value = "2.364.66"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'110.33"
$ws.Range("E5").Value = "  +3.67%  "
$ws.Range("D6").Value = "'308.98"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("D7").Value = "'0.630"
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D9").Value = "'0.617"
$ws.Range("E9").Value = "  +1.60%  "
$ws.Range("D10").Value = "'41.38"
$ws.Range("E10").Value = "  +2.99%  "
$ws.Range("D11").Value = "'0.0920"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "'8.52"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("E13").Value = "  +1.29%  "
$ws.Range("D14").Value = "'0.987"
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").Value = "2.725.11"
$ws.Range("E15").Value = "  +2.33%  "
$ws.Range("D16").Value = "'15.42"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "2.373.14"
$ws.Range("E17").Value = "  +2.39%  "
$ws.Range("D18").Value = "45.392.12"
$ws.Range("E18").Value = "  +5.60%  "
$ws.Range("D19").Value = "'7.31"
$ws.Range("E19").Value = "  -2.16%  "
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").Value = "'13.13"
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("D22").Value = "'73.57"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("E23").Value = "  -1.08%  "
$ws.Range("D24").Value = "'260.27"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("E25").Value = "  +1.80%  "
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("D28").Value = "'7.37"
$ws.Range("E28").Value = "  -4.70%  "
$ws.Range("D29").Value = "'2.36"
$ws.Range("E29").Value = "  +2.33%  "
$ws.Range("D30").Value = "'22.45"
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("E31").Value = "  +10.98%  "
$ws.Range("D32").Value = "'37.96"
$ws.Range("E32").Value = "  -2.02%  "
$ws.Range("D33").Value = "'170.16"
$ws.Range("E33").Value = "  +2.68%  "
$ws.Range("D34").Value = "'2.95"
$ws.Range("E34").Value = "  +7.44%  "
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("E36").Value = "  +3.98%  "
$ws.Range("E37").Value = "  +3.85%  "
$ws.Range("D38").Value = "'3.00"
$ws.Range("E38").Value = "  +7.36%  "
$ws.Range("D39").Value = "'0.0358"
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").Value = "'3.90"
$ws.Range("E40").Value = "  +7.30%  "
$ws.Range("D41").Value = "'1.72"
$ws.Range("E41").Value = "  +8.61%  "
$ws.Range("D42").Value = "'100.72"
$ws.Range("E42").Value = "  -3.04%  "
$ws.Range("D43").Value = "'0.233"
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("D44").Value = "'70.05"
$ws.Range("E44").Value = "  -1.50%  "
$ws.Range("D45").Value = "'13.02"
$ws.Range("E45").Value = "  +6.21%  "
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").Value = "'82.13"
$ws.Range("E47").Value = "  +7.35%  "
$ws.Range("D48").Value = "'9.43"
$ws.Range("E48").Value = "  +6.23%  "
$ws.Range("D49").Value = "'113.11"
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("D50").Value = "'5.56"
$ws.Range("E50").Value = "  +6.46%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.628.65"
$ws.Range("E51").Value = "  -3.83%  "
